$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 115: Primera quality Locoto record
$ws.Cells.Item(115, 1).Value = 1
$ws.Cells.Item(115, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(115, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(115, 4).Value = 44753
$ws.Cells.Item(115, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(115, 5).Value = 15
$ws.Cells.Item(115, 6).Value = 100112042
$ws.Cells.Item(115, 7).Value = "Locoto"
$ws.Cells.Item(115, 8).Value = "Sin especificar"
$ws.Cells.Item(115, 9).Value = "Primera"
$ws.Cells.Item(115, 10).Value = 130
$ws.Cells.Item(115, 11).Value = 37000
$ws.Cells.Item(115, 12).Value = 38000
$ws.Cells.Item(115, 13).Value = 37500
$ws.Cells.Item(115, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(115, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(115, 16).Value = 1875
$ws.Cells.Item(115, 17).Value = 20
$ws.Cells.Item(115, 18).Value = "Hortaliza"

# Row 116: Segunda quality Locoto record
$ws.Cells.Item(116, 1).Value = 1
$ws.Cells.Item(116, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(116, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(116, 4).Value = 44753
$ws.Cells.Item(116, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(116, 5).Value = 15
$ws.Cells.Item(116, 6).Value = 100112042
$ws.Cells.Item(116, 7).Value = "Locoto"
$ws.Cells.Item(116, 8).Value = "Sin especificar"
$ws.Cells.Item(116, 9).Value = "Segunda"
$ws.Cells.Item(116, 10).Value = 200
$ws.Cells.Item(116, 11).Value = 33000
$ws.Cells.Item(116, 12).Value = 35000
$ws.Cells.Item(116, 13).Value = 34000
$ws.Cells.Item(116, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(116, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(116, 16).Value = 1700
$ws.Cells.Item(116, 17).Value = 20
$ws.Cells.Item(116, 18).Value = "Hortaliza"
